$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 18 with delete-team action data
$ws.Range("A18").Value = "2025-08-20 07:58:27"
$ws.Range("B18").Value = "delete-team"
$ws.Range("C18").Value = "new-organization97"
$ws.Range("D18").Value = "newteam"
$ws.Range("E18").Value = "demo"
$ws.Range("I18").Value = "'False"
$ws.Range("I18").Style = "Normal"
